# Scheduled runner update: refresh market-price-derived columns
# (currentAveragePrice / NQ / HQ / LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ)
# for a batch of leve rows across multiple crafting-job sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1817.7407
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 1817.7407
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 5453.2221
$ws.Range("M17").ClearContents()
$ws.Range("N17").Value = -5789.2221

$ws.Range("H97").Value = 16676333
$ws.Range("J97").Value = 16676333
$ws.Range("L97").Value = 50028999
$ws.Range("N97").Value = -50029991

$ws.Range("H100").Value = 1793.3334
$ws.Range("I100").Value = 1270.75
$ws.Range("J100").Value = 2838.5
$ws.Range("K100").Value = 1270.75
$ws.Range("L100").Value = 2838.5
$ws.Range("M100").Value = -729.75
$ws.Range("N100").Value = -3920.5

$ws.Range("H106").Value = 38463436
$ws.Range("I106").Value = 55557196
$ws.Range("J106").Value = 2479.875
$ws.Range("K106").Value = 55557196
$ws.Range("L106").Value = 2479.875
$ws.Range("M106").Value = -55556565
$ws.Range("N106").Value = -3741.875

$ws.Range("H112").Value = 3752.26
$ws.Range("J112").Value = 3752.26
$ws.Range("L112").Value = 11256.78
$ws.Range("N112").Value = -13472.78

$ws.Range("H137").Value = 5241.8623
$ws.Range("I137").Value = 3875.5417
$ws.Range("J137").Value = 11800.2
$ws.Range("K137").Value = 11626.6251
$ws.Range("L137").Value = 35400.60000000001
$ws.Range("M137").Value = -9076.625100000001
$ws.Range("N137").Value = -40500.60000000001

$ws.Range("H138").Value = 3472.721
$ws.Range("I138").Value = 1364.7435
$ws.Range("J138").Value = 5221.8936
$ws.Range("K138").Value = 4094.2305
$ws.Range("L138").Value = 15665.6808
$ws.Range("M138").Value = 1045.7695
$ws.Range("N138").Value = -25945.6808

$ws.Range("H141").Value = 1228.6129
$ws.Range("I141").Value = 370.95834
$ws.Range("J141").Value = 4169.143
$ws.Range("K141").Value = 1112.87502
$ws.Range("L141").Value = 12507.429
$ws.Range("M141").Value = 4067.12498
$ws.Range("N141").Value = -22867.429

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4438.9385
$ws.Range("I32").Value = 3726.8281
$ws.Range("K32").Value = 3726.8281
$ws.Range("M32").Value = -3439.8281

$ws.Range("H45").Value = 3065.25
$ws.Range("I45").Value = 1500.2
$ws.Range("J45").Value = 5673.6665
$ws.Range("K45").Value = 1500.2
$ws.Range("L45").Value = 5673.6665
$ws.Range("M45").Value = -1123.2
$ws.Range("N45").Value = -6427.6665

$ws.Range("H61").Value = 50008810
$ws.Range("I61").Value = 2130.889
$ws.Range("J61").Value = 90923370
$ws.Range("K61").Value = 2130.889
$ws.Range("L61").Value = 90923370
$ws.Range("M61").Value = -1918.889
$ws.Range("N61").Value = -90923794

$ws.Range("H63").Value = 2576.75
$ws.Range("I63").Value = 3300.3333
$ws.Range("J63").Value = 406
$ws.Range("K63").Value = 3300.3333
$ws.Range("L63").Value = 406
$ws.Range("M63").Value = -2614.3333
$ws.Range("N63").Value = -1778

$ws.Range("H66").Value = 2576.75
$ws.Range("I66").Value = 3300.3333
$ws.Range("J66").Value = 406
$ws.Range("K66").Value = 16501.6665
$ws.Range("L66").Value = 2030
$ws.Range("M66").Value = -13069.6665
$ws.Range("N66").Value = -8894

$ws.Range("H136").Value = 50008810
$ws.Range("I136").Value = 2130.889
$ws.Range("J136").Value = 90923370
$ws.Range("K136").Value = 6392.667
$ws.Range("L136").Value = 272770110
$ws.Range("M136").Value = -3842.667
$ws.Range("N136").Value = -272775210

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 8936455
$ws.Range("I134").Value = 22729366
$ws.Range("K134").Value = 68188098
$ws.Range("M134").Value = -68185563

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 10428.071
$ws.Range("J99").Value = 8399.4
$ws.Range("L99").Value = 8399.4
$ws.Range("N99").Value = -11395.4

$ws.Range("H122").Value = 2922.0476
$ws.Range("I122").Value = 1913.7333
$ws.Range("J122").Value = 5442.8335
$ws.Range("K122").Value = 5741.199900000001
$ws.Range("L122").Value = 16328.5005
$ws.Range("M122").Value = -3291.199900000001
$ws.Range("N122").Value = -21228.5005

$ws.Range("H126").Value = 10428.071
$ws.Range("J126").Value = 8399.4
$ws.Range("L126").Value = 25198.2
$ws.Range("N126").Value = -30138.2

$ws.Range("H134").Value = 5091.5
$ws.Range("I134").Value = 2075.0571
$ws.Range("J134").Value = 10118.904
$ws.Range("K134").Value = 6225.1713
$ws.Range("L134").Value = 30356.712
$ws.Range("M134").Value = -3690.1713
$ws.Range("N134").Value = -35426.712

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 9272.125
$ws.Range("I113").Value = 1392.3334
$ws.Range("J113").Value = 14000
$ws.Range("K113").Value = 4177.0002
$ws.Range("L113").Value = 42000
$ws.Range("M113").Value = -2007.0002
$ws.Range("N113").Value = -46340

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 2397.25
$ws.Range("I97").Value = 1814.5
$ws.Range("J97").Value = 2980
$ws.Range("K97").Value = 1814.5
$ws.Range("L97").Value = 2980
$ws.Range("M97").Value = -1318.5
$ws.Range("N97").Value = -3972

$ws.Range("H132").Value = 4266.5938
$ws.Range("I132").Value = 1462.7693
$ws.Range("J132").Value = 16416.5
$ws.Range("K132").Value = 4388.3079
$ws.Range("L132").Value = 49249.5
$ws.Range("M132").Value = -1858.3079
$ws.Range("N132").Value = -54309.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 5141.1
$ws.Range("I100").Value = 2722
$ws.Range("J100").Value = 5745.875
$ws.Range("K100").Value = 2722
$ws.Range("L100").Value = 5745.875
$ws.Range("M100").Value = -2181
$ws.Range("N100").Value = -6827.875

$ws.Range("H132").Value = 7697692.5
$ws.Range("I132").Value = 15154267
$ws.Range("J132").Value = 8099.8438
$ws.Range("K132").Value = 45462801
$ws.Range("L132").Value = 24299.5314
$ws.Range("M132").Value = -45460271
$ws.Range("N132").Value = -29359.5314

$ws.Range("H135").Value = 99999
$ws.Range("J135").Value = 99999
$ws.Range("L135").Value = 99999
$ws.Range("N135").Value = -110139

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 19609444
$ws.Range("I107").Value = 1166.6666
$ws.Range("K107").Value = 3499.9998
$ws.Range("M107").Value = -1579.9998

$ws.Range("H122").Value = 271324.4
$ws.Range("I122").Value = 367360.53
$ws.Range("J122").Value = 7225
$ws.Range("K122").Value = 1102081.59
$ws.Range("L122").Value = 21675
$ws.Range("M122").Value = -1099631.59
$ws.Range("N122").Value = -26575

$ws.Range("H136").Value = 16836298
$ws.Range("I136").Value = 24391226
$ws.Range("K136").Value = 73173678
$ws.Range("M136").Value = -73171128

